$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new log rows (27, 28) describing a finished "consulta" module run.
# Populate column by column to mirror how the original tool wrote the data.

# Column A - ID
$ws.Cells.Item(27,1).Value2 = "07b85aa5-56ed-404b-9827-775c03b79eab"
$ws.Cells.Item(28,1).Value2 = "e17a169a-24cc-4dd7-86dc-fdad5b9f7033"

# Column B - Data/Hora
$ws.Cells.Item(27,2).Value2 = "10/05/2024 20:59:12"
$ws.Cells.Item(28,2).Value2 = "10/05/2024 20:59:20"

# Column C - Computador
$ws.Cells.Item(27,3).Value2 = "DF6963612DTN"
$ws.Cells.Item(28,3).Value2 = "DF6963612DTN"

# Column D - Usuário
$ws.Cells.Item(27,4).Value2 = "anatel_master"
$ws.Cells.Item(28,4).Value2 = "anatel_master"

# Column E - Homologação
$ws.Cells.Item(27,5).Value2 = "01436-19-00953"
$ws.Cells.Item(28,5).Value2 = "01436-19-00953"

# Column F - Atributo
$ws.Cells.Item(27,6).Value2 = "Fornecedor"
$ws.Cells.Item(28,6).Value2 = "WordCloud"

# Column G - Valor
$ws.Cells.Item(27,7).Value2 = "Casas Bahia"
$ws.Cells.Item(28,7).Value2 = "Samsung: 14, A20: 11, DS: 11, Galaxy: 9, SM-A205G: 7, a205g: 7, 32GB: 6, Past: 5, specs: 5, Sm: 4, Battery: 3, price: 3, 4000mAh: 2, Dual: 2, Especificaciones: 2, HD: 2, Images: 2, PhoneMore: 2, Red: 2, SIM: 2, View: 2, eBay: 2, galaxy-a20: 2, hour: 2, itm: 2"

# Column H - Situação
$ws.Cells.Item(27,8).Value2 = 1
$ws.Cells.Item(28,8).Value2 = 1
